$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Marca" (brand) values: Gatito -> Gato, Cachorrito -> Cachorro
# Rows 2-11 hold the "Gatito" product line, rows 12-21 the "Cachorrito" line.
For ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value = "Gato"
}
For ($r = 12; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = "Cachorro"
}

# Widen the "Nome produto" column (B) to fit the text.
$ws.Range("B:B").ColumnWidth = 21.14

# Stray formatting leftover near the bottom of the sheet (underlined cell D23).
$ws.Range("D23").Font.Underline = 2
$ws.Rows(23).RowHeight = 15.75

# Leave the cursor parked on B22, matching the saved selection state.
$ws.Range("B22").Select()
